$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the FV2210 / FV2304 header labels (formerly *_old / *_new) ---
$headers = @(
  "Segmentname_FV2210","Segmentgruppe_FV2210","Segment_FV2210","Datenelement_FV2210","Segment ID_FV2210",
  "Code_FV2210","Qualifier_FV2210","Beschreibung_FV2210","Bedingungsausdruck_FV2210","Bedingung_FV2210",
  "diff",
  "Segmentname_FV2304","Segmentgruppe_FV2304","Segment_FV2304","Datenelement_FV2304","Segment ID_FV2304",
  "Code_FV2304","Qualifier_FV2304","Beschreibung_FV2304","Bedingungsausdruck_FV2304","Bedingung_FV2304"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Freeze the header row ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the data range into an Excel Table, preserving the existing  ---
#        header-row formatting instead of letting a fresh dxf get captured.
$headerRange = $ws.Range("A1:U1")
$scratchRange = $ws.Range("A100:U100")
$headerRange.Copy($scratchRange)
$headerRange.ClearFormats()

$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U59"), 0, 1)
$tbl.Name = "Table1"

$scratchRange.Copy()
$headerRange.PasteSpecial(-4122)
$scratchRange.Clear()
$ws.Range("A1").Select()

$tbl.TableStyle = ""
